$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Description") is empty in every data row; remove it entirely.
# This shifts the Receipt column (and its hyperlinks) from F to E.

# Remember each existing hyperlink's row + target URL before we touch
# anything, since the column delete does not re-home the hyperlink ranges.
$links = @()
foreach ($hl in $ws.Hyperlinks) {
    $links += [PSCustomObject]@{ Row = $hl.Range.Row; Target = $hl.Address }
}

$ws.Columns.Item(5).Delete()

# Drop the now stale hyperlinks (still anchored on the old column F) and
# recreate them on column E, restoring the Hyperlink cell style afterwards
# since Hyperlinks.Add() mints a fresh (slightly different) style entry.
$ws.Hyperlinks.Delete()
foreach ($l in $links) {
    $cell = $ws.Cells.Item($l.Row, 5)
    $ws.Hyperlinks.Add($cell, $l.Target) | Out-Null
    $cell.Style = "Hyperlink"
}
